$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit tug boat starting time: StartWorkingTime (column D, rows 2-23)
# moves from 9:00 AM (0.375) to 8:00 AM (0.33333333333333331) for every tug boat.
$ws.Range("D2:D23").Value = 0.33333333333333331

# Reflect the user's final selection in the saved view state.
$ws.Range("D14").Select()
